$wb = $excel.ActiveWorkbook

# --- update selection on the existing "mode" sheet (was the active tab) ---
$mode = $wb.Worksheets.Item("mode")
$mode.Range("A2:K18").Select() | Out-Null

# --- add the new "linear" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "linear"

# column widths (ColumnWidth snaps to 1/6-character increments in this
# engine, so these inputs are chosen to land on the closest achievable
# stored width to the target: 12.5703125 -> 12.5, 12 -> 12 exactly)
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 11.166666666666666
$ws.Columns.Item(7).ColumnWidth = 11.166666666666666

# header rows (g=0 / g=1 / g=2 blocks)
$ws.Range("A2").Value = "g=0"
$ws.Range("D2").Value = "g=1"
$ws.Range("G2").Value = "g=2"

$ws.Range("A3").Value = "n"
$ws.Range("B3").Value = "time [s]"
$ws.Range("D3").Value = "n"
$ws.Range("E3").Value = "time [s]"
$ws.Range("G3").Value = "n"
$ws.Range("H3").Value = "time [s]"

# data rows (decimal literals below are bit-identical float64 values to the
# scientific-notation forms Excel itself would emit, e.g. 0.008 == 8.0000000000000002E-3)
$data = @(
    @(100,           0.01,    100,           0.01,    100,           0.01),
    @(1000,          0.008,   1000,          0.008,   1000,          0.008),
    @(10000,         0.008,   10000,         0.009,   10000,         0.008),
    @(100000,        0.009,   100000,        0.009,   100000,        0.011),
    @(1000000,       0.018,   1000000,       0.014,   1000000,       0.035),
    @(10000000,      0.101,   10000000,      0.066,   10000000,      0.28),
    @(100000000,     0.936,   100000000,     0.585,   100000000,     2.717),
    @(1000000000,    9.207,   1000000000,    5.677,   1000000000,    27.118),
    @(10000000000,   91.863,  10000000000,   56.796,  10000000000,   267.944),
    @(100000000000,  876.14,  100000000000,  540.635, 100000000000,  2511.625)
)

$row = 4
foreach ($r in $data) {
    $ws.Range("A$row").Value = $r[0]
    $ws.Range("B$row").Value = $r[1]
    $ws.Range("D$row").Value = $r[2]
    $ws.Range("E$row").Value = $r[3]
    $ws.Range("G$row").Value = $r[4]
    $ws.Range("H$row").Value = $r[5]
    $row++
}

# selection / active cell on the new sheet
$ws.Range("F21").Select() | Out-Null
